# Edit script: apply proofErr markers + paragraph-8 rewrite + new Untrained/Trained Weights section
$d = $word.ActiveDocument

# --- 1) Wrap "Siska" in the "Professor Siska" paragraph with spellStart/spellEnd proofErr markers ---
$siskaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.TrimEnd() -eq "Professor Siska") {
        $siskaPara = $candidate
        break
    }
}
if ($siskaPara -ne $null) {
    $siskaXml = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Professor</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Siska</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
    $siskaPara.Range.InsertXML($siskaXml)
}

# --- 2) Rewrite the reflection paragraph (javathcript / api / cdr / Javathcript / far proofErr splits)
#        and append the new "Untrained Weights" / weight-dump / "Trained Weights" paragraphs after it ---
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("I used the professor")) {
        $targetPara = $candidate
        break
    }
}
if ($targetPara -ne $null) {
    $replacementXml = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">I used the professor’s skeleton code. I tried to use the book’s pseudocode on A* as reference. The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>javathcript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> library has bare emacs lisp </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> like car, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>cdr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, and list.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Javathcript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> did not have the sort function I needed to sort the open list.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> So </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>far</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the explorer has reached the destination.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Untrained Weights</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>self.bias</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = [list([0.3, 0.5, 0.6, 0.7]),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                     list([0.2, 0.1, 0.3, 0.2, 0.4, 0.5, 0.6, 0.2])]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>self.weights</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = [list([[.3, -.45, 0.12, -0.12, 0.89, -0.25, 0.16, -0.56, 0.12, 0.43],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.12, -0.23, -0.15, 0.62, 0.821, -</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                                  0.512, 0.123, 0.321, -0.9, .5],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [-0.64, 0.52, 0.24, 0.93, -0.84, -</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                                  0.64, -0.52, 0.34, -0.6, -0.7],</w:t></w:r></w:p><w:p><w:lastRenderedPageBreak/><w:r><w:t xml:space="preserve">                              [-0.52, 0.12, 0.84, -0.51, 0.321, -0.123, -0.721, 0.632, 0.5, 0.3]]),</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                        list([[0.123, -0.21, -0.21, 0.123],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.412, 0.31, -0.731, -0.365],    </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.231, 0.521, -0.441, -0.3852],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.521, -0.41, -0.581, 0.112],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.412, -0.891, 0.712, -0.321],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.221, 0.551, -0.213, -0.555],</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                              [0.812, -0.231, 0.367, -0.888],</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">                              [-0.512, -0.231, 0.751, 0.121]])]</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Trained Weights</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>'
    $targetPara.Range.InsertXML($replacementXml)
}
